# Apply "Updated symbol list" edits to the crypto price sheet.
# All touched cells are plain text (inline strings) in columns D (Price)
# and E (Volume(1h)) of Sheet1. The Price column values look numeric, so
# they must be forced to text storage, matching the source file, and then
# restored to the default (unstyled) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue "D2"  "245.63"
Set-TextValue "D3"  "25.25"
Set-TextValue "D5"  "0.05598"
Set-TextValue "D6"  "6.558"
Set-TextValue "D8"  "0.8150"
Set-TextValue "D9"  "0.8352"
Set-TextValue "D11" "0.06958"
Set-TextValue "D13" "0.02829"
Set-TextValue "D14" "0.09389"
Set-TextValue "D16" "0.0005950"
Set-TextValue "D17" "0.006136"
Set-TextValue "D18" "3.497"
Set-TextValue "D22" "3.750"
Set-TextValue "D23" "0.04683"
Set-TextValue "D25" "0.001243"
Set-TextValue "D40" "0.03666"
Set-TextValue "D42" "0.1058"
Set-TextValue "D44" "0.008223"
Set-TextValue "D45" "0.00005293"

# Column E (Volume(1h)) updates
$ws.Range("E16").Value = "15OneONEWorstin24h"
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
